# Generate Report for Handoff
#
# - Updates the "in progress" status/timestamp cells to reflect a completed
#   handoff ("Ready for handoff").
# - Removes the second data row (the ed864e02-... file) from every sheet,
#   including its now-orphaned hyperlinks.

function Remove-RowHyperlinks {
    param($ws, [int]$row)

    $again = $true
    while ($again) {
        $again = $false
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row) {
                $h.Delete()
                $again = $true
                break
            }
        }
    }
}

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-41-18 20:41:50"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-18 20:41:47"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-18 20:41:50"

# --- Remove the ed864e02-... row (row 3) from every sheet -------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Rows.Item(3).Delete()
    Remove-RowHyperlinks $ws 3
}
